# Insert a new "汽車" (car) worksheet between "建物" (building) and "保險" (insurance),
# matching the commit that split a combined sheet into land / building / car / insurance / debt.

$wb = $excel.ActiveWorkbook

# Locate the "建物" sheet so the new sheet lands right after it (and therefore right
# before "保險"), producing the order: 土地, 建物, 汽車, 保險, 債務.
$buildingSheet = $wb.Worksheets.Item("建物")
$carSheet = $wb.Worksheets.Add($null, $buildingSheet)
$carSheet.Name = "汽車"

# ---- Header row (row 1), styled like the other sheets: bold, centered, top-aligned, bordered ----
# Copy the formatting straight from the "建物" header row / index column instead of
# re-building it property-by-property, so we reuse the workbook's existing style
# entries instead of minting new (redundant) ones.
$buildingSheet.Range("B1:N1").Copy() | Out-Null
$carSheet.Range("B1:N1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$buildingSheet.Range("A2").Copy() | Out-Null
$carSheet.Range("A2").PasteSpecial(-4122) | Out-Null      # xlPasteFormats

$carSheet.Range("B1").Value = "name"
$carSheet.Range("C1").Value = "capacity"
$carSheet.Range("D1").Value = "owner"
$carSheet.Range("E1").Value = "register_date"
$carSheet.Range("F1").Value = "register_reason"
$carSheet.Range("G1").Value = "acquire_value"
$carSheet.Range("H1").Value = "property_category"
$carSheet.Range("I1").Value = "category"
$carSheet.Range("J1").Value = "date"
$carSheet.Range("K1").Value = "legislator_name"
$carSheet.Range("L1").Value = "legislator_id"
$carSheet.Range("M1").Value = "source_file"
$carSheet.Range("N1").Value = "index"

# ---- Data row (row 2) ----
$carSheet.Range("A2").Value = 31
$carSheet.Range("B2").Value = "國瑞(國產091603)"
$carSheet.Range("C2").Value = 1497
$carSheet.Range("D2").Value = "何欣純"
$carSheet.Range("E2").Value = "96年10月25日"
$carSheet.Range("F2").Value = "買賣"
$carSheet.Range("G2").Value = "(超過五年）"
$carSheet.Range("H2").Value = "car"
$carSheet.Range("I2").Value = "normal"

# Force the "date" column to remain plain text "2013-12-30" instead of being
# auto-converted into a date serial number by Excel's smart entry.
$carSheet.Range("J2").NumberFormat = "@"
$carSheet.Range("J2").Value = "2013-12-30"

$carSheet.Range("K2").Value = "何欣純"
$carSheet.Range("L2").Value = 1733
$carSheet.Range("M2").Value = "tmp8e3c1"
$carSheet.Range("N2").Value = 31
